$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 7.143069986577376
$ws.Range("D2").Value = 9.159184821805955
$ws.Range("E2").Value = 12.2588898689187
$ws.Range("F2").Value = 45.91109668700077
$ws.Range("G2").Value = 3.745819127465058
$ws.Range("I2").Value = 37.36315893148684
$ws.Range("L2").Value = 9.172474489855926
$ws.Range("M2").Value = 61.50499109747015
$ws.Range("C3").Value = 6.785600552501571
$ws.Range("D3").Value = 9.396766635983203
$ws.Range("E3").Value = 11.74367124332952
$ws.Range("F3").Value = 46.42636363892459
$ws.Range("G3").Value = 3.756711958689575
$ws.Range("I3").Value = 37.99083152135487
$ws.Range("L3").Value = 8.988510766419925
$ws.Range("M3").Value = 58.4494859549805
$ws.Range("C4").Value = 6.556065328879683
$ws.Range("D4").Value = 9.548547313203844
$ws.Range("E4").Value = 11.41789697102617
$ws.Range("F4").Value = 46.78740219708561
$ws.Range("G4").Value = 3.763655207859551
$ws.Range("I4").Value = 38.40779295640598
$ws.Range("L4").Value = 8.876638719217441
$ws.Range("M4").Value = 56.49051309512775
$ws.Range("C5").Value = 6.460065807844111
$ws.Range("D5").Value = 9.611880233223214
$ws.Range("E5").Value = 11.28293393676064
$ws.Range("F5").Value = 46.94535431086861
$ws.Range("G5").Value = 3.766549638614326
$ws.Range("I5").Value = 38.58540398967893
$ws.Range("L5").Value = 8.831369881879308
$ws.Range("M5").Value = 55.67196471410378
$ws.Range("C6").Value = 6.44397860161051
$ws.Range("D6").Value = 9.622485995191541
$ws.Range("E6").Value = 11.26039533385816
$ws.Range("F6").Value = 46.97222443372878
$ws.Range("G6").Value = 3.767034207582197
$ws.Range("I6").Value = 38.61535422699233
$ws.Range("L6").Value = 8.823873764880906
$ws.Range("M6").Value = 55.53484190075149
$ws.Range("C7").Value = 6.554780519761459
$ws.Range("D7").Value = 9.549395449842885
$ws.Range("E7").Value = 11.41608551716916
$ws.Range("F7").Value = 46.78948908511659
$ws.Range("G7").Value = 3.763693978840938
$ws.Range("I7").Value = 38.41015742787023
$ws.Range("L7").Value = 8.876026848566484
$ws.Range("M7").Value = 56.47955499291582
$ws.Range("C8").Value = 7.021934654534914
$ws.Range("D8").Value = 9.239874791004702
$ws.Range("E8").Value = 12.0832902321862
$ws.Range("F8").Value = 46.07924330938925
$ws.Range("G8").Value = 3.749522587308725
$ws.Range("I8").Value = 37.57287740520184
$ws.Range("L8").Value = 9.108847385697102
$ws.Range("M8").Value = 60.4689722240232
$ws.Range("C9").Value = 7.856084700352318
$ws.Range("D9").Value = 8.680022038704283
$ws.Range("E9").Value = 13.31083146439406
$ws.Range("F9").Value = 45.05961640663357
$ws.Range("G9").Value = 3.723714437701076
$ws.Range("I9").Value = 36.19324851681284
$ws.Range("L9").Value = 9.571982773317183
$ws.Range("M9").Value = 67.61442191613081
$ws.Range("C10").Value = 8.416755245016303
$ws.Range("D10").Value = 8.297994235674299
$ws.Range("E10").Value = 14.15680581396617
$ws.Range("F10").Value = 44.56461234069673
$ws.Range("G10").Value = 3.70590097007736
$ws.Range("I10").Value = 35.35727319514808
$ws.Range("L10").Value = 9.91372579515415
$ws.Range("M10").Value = 72.43089735334578
$ws.Range("C11").Value = 8.660192319454829
$ws.Range("D11").Value = 8.130767798681012
$ws.Range("E11").Value = 14.52850146991964
$ws.Range("F11").Value = 44.4005432652246
$ws.Range("G11").Value = 3.698032451873747
$ws.Range("I11").Value = 35.01986245376937
$ws.Range("L11").Value = 10.06903485049191
$ws.Range("M11").Value = 74.52535540453763
$ws.Range("C12").Value = 8.750689444895015
$ws.Range("D12").Value = 8.068411977783137
$ws.Range("E12").Value = 14.66729916188889
$ws.Range("F12").Value = 44.34770271162694
$ws.Range("G12").Value = 3.695085436799638
$ws.Range("I12").Value = 34.89865440750911
$ws.Range("L12").Value = 10.12778878455516
$ws.Range("M12").Value = 75.30445081375127
$ws.Range("C13").Value = 8.731274518600838
$ws.Range("D13").Value = 8.081797828138344
$ws.Range("E13").Value = 14.63749445347214
$ws.Range("F13").Value = 44.35866171557831
$ws.Range("G13").Value = 3.695718697440979
$ws.Range("I13").Value = 34.92446026258752
$ws.Range("L13").Value = 10.11513820395294
$ws.Range("M13").Value = 75.13728435972179
$ws.Range("C14").Value = 8.667671463569334
$ws.Range("D14").Value = 8.12561815084573
$ws.Range("E14").Value = 14.53995995287697
$ws.Range("F14").Value = 44.39600699066843
$ws.Range("G14").Value = 3.697789352225681
$ws.Range("I14").Value = 35.00975670905084
$ws.Range("L14").Value = 10.07386993455052
$ws.Range("M14").Value = 74.58973386008483
$ws.Range("C15").Value = 8.628492625336218
$ws.Range("D15").Value = 8.152586466849348
$ws.Range("E15").Value = 14.47996090648242
$ws.Range("F15").Value = 44.4201066923038
$ws.Range("G15").Value = 3.699061901016973
$ws.Range("I15").Value = 35.06287007223602
$ws.Range("L15").Value = 10.04858332620636
$ws.Range("M15").Value = 74.25251292845935
$ws.Range("C16").Value = 8.40061075045671
$ws.Range("D16").Value = 8.309056941722408
$ws.Range("E16").Value = 14.13224423685608
$ws.Range("F16").Value = 44.57660706391708
$ws.Range("G16").Value = 3.706419825658047
$ws.Range("I16").Value = 35.38021982035157
$ws.Range("L16").Value = 9.90356973410154
$ws.Range("M16").Value = 72.29206269107237
$ws.Range("C17").Value = 8.257823759709138
$ws.Range("D17").Value = 8.406743041950534
$ws.Range("E17").Value = 13.91551395746665
$ws.Range("F17").Value = 44.68862165584409
$ws.Range("G17").Value = 3.710993035870332
$ws.Range("I17").Value = 35.5861586862612
$ws.Range("L17").Value = 9.814542686235475
$ws.Range("M17").Value = 71.06453387736505
$ws.Range("C18").Value = 8.17460383116058
$ws.Range("D18").Value = 8.463545040249784
$ws.Range("E18").Value = 13.78962235364948
$ws.Range("F18").Value = 44.7587710827291
$ws.Range("G18").Value = 3.713645609745357
$ws.Range("I18").Value = 35.70861975354595
$ws.Range("L18").Value = 9.763324373642961
$ws.Range("M18").Value = 70.34940573608576
$ws.Range("C19").Value = 8.146239949794259
$ws.Range("D19").Value = 8.482882257434518
$ws.Range("E19").Value = 13.74678790210135
$ws.Range("F19").Value = 44.78348998474402
$ws.Range("G19").Value = 3.714547570221098
$ws.Range("I19").Value = 35.75076016707748
$ws.Range("L19").Value = 9.745981789397867
$ws.Range("M19").Value = 70.10572036807875
$ws.Range("C20").Value = 8.27313690343467
$ws.Range("D20").Value = 8.396280317628584
$ws.Range("E20").Value = 13.93871350895271
$ws.Range("F20").Value = 44.67610173449835
$ws.Range("G20").Value = 3.71050392205255
$ws.Range("I20").Value = 35.56381829475119
$ws.Range("L20").Value = 9.824021317164156
$ws.Range("M20").Value = 71.19614799383569
$ws.Range("C21").Value = 8.686399119952942
$ws.Range("D21").Value = 8.112720515666814
$ws.Range("E21").Value = 14.56866171983931
$ws.Range("F21").Value = 44.38478166312886
$ws.Range("G21").Value = 3.697180274900671
$ws.Range("I21").Value = 34.98452168162938
$ws.Range("L21").Value = 10.08599327797638
$ws.Range("M21").Value = 74.75094435369033
$ws.Range("C22").Value = 8.946651543671919
$ws.Range("D22").Value = 7.933064225336517
$ws.Range("E22").Value = 14.9689512324428
$ws.Range("F22").Value = 44.24876489438152
$ws.Range("G22").Value = 3.688662066221125
$ws.Range("I22").Value = 34.64437111759219
$ws.Range("L22").Value = 10.34604843135982
$ws.Range("M22").Value = 76.99240114427297
$ws.Range("C23").Value = 8.808654294956805
$ws.Range("D23").Value = 8.028421288519979
$ws.Range("E23").Value = 14.75637194938778
$ws.Range("F23").Value = 44.31621735489527
$ws.Range("G23").Value = 3.693191447836799
$ws.Range("I23").Value = 34.82225899239202
$ws.Range("L23").Value = 10.1791126062029
$ws.Range("M23").Value = 75.80361186305906
$ws.Range("C24").Value = 8.266217348832672
$ws.Range("D24").Value = 8.401008520718982
$ws.Range("E24").Value = 13.92822900615205
$ws.Range("F24").Value = 44.68174411657297
$ws.Range("G24").Value = 3.710724977516437
$ws.Range("I24").Value = 35.57390575672028
$ws.Range("L24").Value = 9.819736135814129
$ws.Range("M24").Value = 71.13667453634366
$ws.Range("C25").Value = 7.63946244575946
$ws.Range("D25").Value = 8.826386944756139
$ws.Range("E25").Value = 12.98816160793515
$ws.Range("F25").Value = 45.29280156901859
$ws.Range("G25").Value = 3.730490131309864
$ws.Range("I25").Value = 36.53665376892417
$ws.Range("L25").Value = 9.446259455622652
$ws.Range("M25").Value = 65.75635686212388
